$wb = $excel.ActiveWorkbook

# --- Update scraped_at timestamps on the "snapshot" sheet (column K, rows 2-34) ---
$snapshot = $wb.Worksheets.Item("snapshot")

$timestamps = @{
    2  = "2025-12-05T07:02:35.041784+00:00"
    3  = "2025-12-05T07:02:35.041823+00:00"
    4  = "2025-12-05T07:02:35.041845+00:00"
    5  = "2025-12-05T07:02:37.297198+00:00"
    6  = "2025-12-05T07:02:37.297225+00:00"
    7  = "2025-12-05T07:02:40.025932+00:00"
    8  = "2025-12-05T07:02:42.354517+00:00"
    9  = "2025-12-05T07:02:44.737842+00:00"
    10 = "2025-12-05T07:02:47.506568+00:00"
    11 = "2025-12-05T07:02:52.530432+00:00"
    12 = "2025-12-05T07:02:52.530461+00:00"
    13 = "2025-12-05T07:02:55.467434+00:00"
    14 = "2025-12-05T07:02:57.777160+00:00"
    15 = "2025-12-05T07:02:57.777190+00:00"
    16 = "2025-12-05T07:03:00.635564+00:00"
    17 = "2025-12-05T07:03:03.410118+00:00"
    18 = "2025-12-05T07:03:03.410146+00:00"
    19 = "2025-12-05T07:03:06.202039+00:00"
    20 = "2025-12-05T07:03:06.202081+00:00"
    21 = "2025-12-05T07:03:08.955609+00:00"
    22 = "2025-12-05T07:03:08.955636+00:00"
    23 = "2025-12-05T07:03:11.305873+00:00"
    24 = "2025-12-05T07:03:11.305903+00:00"
    25 = "2025-12-05T07:03:11.305923+00:00"
    26 = "2025-12-05T07:03:13.609252+00:00"
    27 = "2025-12-05T07:03:18.255218+00:00"
    28 = "2025-12-05T07:03:18.255249+00:00"
    29 = "2025-12-05T07:03:21.020301+00:00"
    30 = "2025-12-05T07:03:21.020331+00:00"
    31 = "2025-12-05T07:03:23.355031+00:00"
    32 = "2025-12-05T07:03:23.355062+00:00"
    33 = "2025-12-05T07:03:26.137939+00:00"
    34 = "2025-12-05T07:03:26.137967+00:00"
}

foreach ($row in $timestamps.Keys) {
    $snapshot.Cells.Item($row, 11).Value = $timestamps[$row]
}

# --- Remove the four newly-injured player rows from "new_injured" sheet, ---
# --- leaving only the header row (A1:G1). ---
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Range("A2:G5").EntireRow.Delete()
